$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the first data row (the "juanito" sample row) - this shifts the
# former "Invoice Total" row (row 3) up to row 2, matching the diff which
# deletes row 2 and drops row 3 entirely (dimension becomes A1:K2).
$ws.Rows.Item(2).Delete()

# The invoice total row now sitting at row 2 should summarize zero
# invoices (Invoice Amount total and 0-30 days bucket both reset to 0).
$ws.Cells.Item(2, 3).Value = 0
$ws.Cells.Item(2, 7).Value = 0
